$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New patient row: Sheldon Cooper
$ws.Range("A3").Value = "Sheldon"
$ws.Range("B3").Value = "Cooper"
$ws.Range("C3").Value = "03/23/1995"
$ws.Range("D3").Value = 3434343434
$ws.Range("E3").Value = "sheldon@bbt.com"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:sheldon@bbt.com")
$ws.Range("F3").Value = "otra dir 324"
$ws.Range("G3").Value = "fx45544555"

# Next row down, just the active/selected cell, left empty
$ws.Range("G4").Value = ""
$ws.Range("G4").Font.Name = "Calibri"

$ws.Range("G4").Select()
